$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "x"
$ws.Range("B3").Value = "x"
$ws.Range("B4").Value = "x"
$ws.Range("B5").Value = "x"
$ws.Range("B6").Value = "x"
